$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New location rows (50-55) -------------------------------------------------
# Columns: A name | B game_map_id | C quest_reward_item_id | D required_quest_item_id
#          E description | F is_port | G can_players_enter | H enemy_strength_type
#          I can_auto_battle | J x | K y | L type | M drops_items | N pin_css_class

# Row 50 - Northren Port
$ws.Range("A50").Value = "Northren Port"
$ws.Range("B50").Value = "Delusional Memories"
$ws.Range("E50").Value = "A port to the north. The men and women here are hardened fighters and vetrans of the battles between the Federation and free people of the south."
$ws.Range("F50").Value = 1
$ws.Range("G50").Value = 1
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = 304
$ws.Range("K50").Value = 192
$ws.Range("M50").Value = "No"

# Row 51 - Southren Port
$ws.Range("A51").Value = "Southren Port"
$ws.Range("B51").Value = "Delusional Memories"
$ws.Range("E51").Value = "A port of the free people in the south. They are known to trade with peoples of the Northren Port for supplies and sending goods and services back and forth. Although these days with The Federation it's harder to trade."
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = 1
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = 304
$ws.Range("K51").Value = 288
$ws.Range("M51").Value = "No"

# Row 52 - Federation Controlled Town
$ws.Range("A52").Value = "Federation Controlled Town"
$ws.Range("B52").Value = "Delusional Memories"
$ws.Range("E52").Value = "The people here are slaves to the Federation. There's the poor and the rich, there is no in between. You either work for The Federation for nothing but scraps, or you are born into a family of one who works for The Federation. Even soldiers sent to die have a higher standing then the people of this town."
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 6
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = 16
$ws.Range("K52").Value = 224
$ws.Range("M52").Value = "Yes"

# Row 53 - Delusional Abandoned Gold Mines
$ws.Range("A53").Value = "Delusional Abandoned Gold Mines"
$ws.Range("B53").Value = "Delusional Memories"
$ws.Range("E53").Value = "These old Gold Mines hold the memories of the past as haunting apperations."
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 6
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 448
$ws.Range("K53").Value = 320
$ws.Range("M53").Value = "Yes"

# Row 54 - Alchemy Corrupted Church
$ws.Range("A54").Value = "Alchemy Corrupted Church"
$ws.Range("B54").Value = "Delusional Memories"
$ws.Range("D54").Value = "Purgatory's Lantern"
$ws.Range("E54").Value = "A church corrupted by the magics of Alchemy. The gates of time have opened here, the judges of time step forth."
$ws.Range("G54").Value = 1
$ws.Range("J54").Value = 400
$ws.Range("K54").Value = 16
$ws.Range("M54").Value = "No"

# Row 55 - Federation City
$ws.Range("A55").Value = "Federation City"
$ws.Range("B55").Value = "Delusional Memories"
$ws.Range("E55").Value = "The main city where the Federation organizes it's military plans from. No army, not even The Red Hawks have managed to breech the city because of the Alchemy that The Church practices, in conjunction with thier religious beliefs."
$ws.Range("G55").Value = 1
$ws.Range("J55").Value = 80
$ws.Range("K55").Value = 96
$ws.Range("M55").Value = "No"

# --- Column width tweaks (A and B got wider to fit the new longer values) -----
$ws.Columns.Item(1).ColumnWidth = 36.833333333333336
$ws.Columns.Item(2).ColumnWidth = 22.666666666666668
